$d = $word.ActiveDocument

# The author fixed a typo ("they way" -> "the way") in the bold prompt
# "How do you think they way you treat your workspace affects others in
# the school?" and, because that's how Word records the last edit
# location, a "_GoBack" bookmark was left right after the corrected word.

# 1) Fix the typo: "they way" -> "the way"
$findRange = $d.Content.Duplicate
$findRange.Find.Execute("they way", $false, $false, $false, $false, $false, `
    $true, 1, $false, "the way", 2)

# 2) Re-locate the corrected phrase so we can drop the bookmark exactly
#    after "...think the" (i.e. where the cursor would have been left).
$markRange = $d.Content.Duplicate
$markRange.Find.Execute("How do you think the", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

$goBack = $d.Range($markRange.End, $markRange.End)
$d.Bookmarks.Add("_GoBack", $goBack)
